# Update the Mouser Part # for the DC power jack (J3) in the BOM.
# The jack was changed from " 163-R123B-E" to "163-1060-EX" because the
# old part had a metal jack body with center-negative wiring, which caused
# a ground short.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

$ws.Range("F12").Value = "163-1060-EX"
